$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF), matching the style of the
# existing header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-32 for columns I and J
$data = @(
    @(2, 1, 5),
    @(3, 1, 5),
    @(4, 1, 6),
    @(5, 1, 6),
    @(6, 5, 10),
    @(7, 1, 7),
    @(8, 1, 7),
    @(9, 1, 7),
    @(10, 1, 5),
    @(11, 1, 6),
    @(12, 1, 3),
    @(13, 1, 4),
    @(14, 5, 8),
    @(15, 8, 8),
    @(16, 7, 8),
    @(17, 6, 6),
    @(18, 4, 5),
    @(19, 8, 8),
    @(20, 4, 6),
    @(21, 5, 6),
    @(22, 4, 8),
    @(23, 7, 7),
    @(24, 6, 6),
    @(25, 3, 5),
    @(26, 2, 3),
    @(27, 8, 8),
    @(28, 9, 9),
    @(29, 4, 5),
    @(30, 7, 7),
    @(31, 4, 5),
    @(32, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
